$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (t_period 2030) gets new scenario probabilities
$ws.Range("B4").Value = 0.62
$ws.Range("C4").Value = 0.3
$ws.Range("D4").Value = 0.08

# Row 5 changes from 2035 to 2040, with the same new probabilities as 2030
$ws.Range("A5").Value = 2040
$ws.Range("B5").Value = 0.62
$ws.Range("C5").Value = 0.3
$ws.Range("D5").Value = 0.08

# New row 6 (t_period 2050) carries the previous 2035 values
$ws.Range("A6").Value = 2050
$ws.Range("B6").Value = 0.81
$ws.Range("C6").Value = 0.15
$ws.Range("D6").Value = 0.04
